# sprint 15 (for pagination and bug fixing)
# Updates the Input_SGST report with the current client / ledger info and
# refreshes the sample invoice row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reports")

# --- simple text fields (no date-like pattern, safe to assign directly) ---
$ws.Range("B2").Value = "rk industries"
$ws.Range("B6").Value = "Monthly"

# --- date-like text fields ---
# Plain `.Value = "2020-09-11"` gets auto-recognised as a real date by the
# engine (like Excel's smart-entry) and turned into a numeric serial with a
# date number format, which is not what the report stores (it keeps these
# as literal text). Forcing NumberFormat="@" (Text) on a scratch cell before
# the assignment keeps it literal text; copying the *value only* from that
# scratch cell over to the real target then leaves the target cell's own
# style/border untouched.
$scratch = $ws.Range("B10")
$scratch.NumberFormat = "@"

$scratch.Value = "2020-09-02"
$scratch.Copy()
$ws.Range("A10").PasteSpecial(-4163)   # xlPasteValues

$scratch.Value = "2020-09-11"
$scratch.Copy()
$ws.Range("B5").PasteSpecial(-4163)    # xlPasteValues

# restore the scratch cell's own formatting (border-only style) before
# clearing it back out, then blank B10/C10 (invoice voucher no./number are
# removed from this row).
$ws.Range("C10").Copy()
$ws.Range("B10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B10").ClearContents()
$ws.Range("C10").ClearContents()

# --- invoice amount ---
$ws.Range("E10").Value = 9.800000000000001
